$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "biaya" column (column B) entirely, shifting subsequent
# columns (keterangan, etc.) one position to the left.
$ws.Columns.Item(2).Delete()

# Move selection to A2 to mirror the saved view state.
$ws.Range("A2").Select()
